# "Data source corrected and updated"
#
# The sheet's J/K columns held a placeholder pair on row 1 ("r"/"s" text,
# via shared strings) and a 0.3/1 numeric pair on every other row. The
# correction replaces all of that with a single consistent numeric pair
# (J=1, K=0.3) down the whole used range, which also makes the
# shared-strings table empty since nothing references text anymore.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 51
For ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 10).Value = 1    # column J
    $ws.Cells.Item($r, 11).Value = 0.3  # column K
}

# Reflect the reviewer's last selection: column K selected, active cell K1,
# scrolled down so row 26 is at the top of the viewport.
$ws.Range("K1:K51").Select()
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 1
